$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "Player Info" worksheet as the first sheet ---------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Header row (ID, NAME, BATTING_HAND, BOWL_STYLE) with bold/bordered header style
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row for player 4739 (ID is numeric-looking, force text formatting so it
# is stored the same way the source data stores every value - as text - then
# reset the number format back to General so no stray formatting is left behind)
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4739"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Carl Junior Dala"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium"

# --- 2. Rename MATCH_CARD_LINK -> MATCH_CODE and replace URLs with match codes --

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4187"
$battingSheet.Range("D2").Style = "Normal"
$battingSheet.Range("D3").NumberFormat = "@"
$battingSheet.Range("D3").Value = "4188"
$battingSheet.Range("D3").Style = "Normal"

$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4187"
$bowlingSheet.Range("B2").Style = "Normal"
$bowlingSheet.Range("B3").NumberFormat = "@"
$bowlingSheet.Range("B3").Value = "4188"
$bowlingSheet.Range("B3").Style = "Normal"

Write-Host "Edit complete"
